$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    ,@(66.6808955223879, 54.09333333333332, 78.34666666666666)
    ,@(1.534726368159203, 0.3600000000000002, 3.320000000000001)
    ,@(6.253432835820892, 2.58, 11.16666666666666)
    ,@(1.784975124378108, 0.4000000000000002, 3.96)
    ,@(25.889983416252, 13.86666666666667, 38.65333333333333)
    ,@(63.14819237147591, 49.16666666666667, 76.11999999999999)
    ,@(89.04019900497501, 79.47333333333336, 94.36000000000004)
    ,@(23.58540630182408, 13.36000000000001, 36.23999999999999)
    ,@(88.90338308457697, 79.76000000000002, 94.33999999999997)
    ,@(77.34553897180754, 66.37333333333335, 85.91333333333333)
    ,@(51.07946932006621, 37.12666666666668, 63.85333333333331)
    ,@(86.99787728026517, 78.09333333333333, 92.93333333333335)
    ,@(83.37290215588713, 73.16000000000001, 89.99333333333334)
    ,@(2.17429519071309, 0.7400000000000003, 4.200000000000002)
    ,@(25.95688225538966, 13.98666666666666, 38.63999999999999)
    ,@(51.13737976782739, 36.63333333333333, 65.80666666666667)
    ,@(6.253432835820892, 2.58, 11.16666666666666)
    ,@(77.95038142620226, 67.46666666666667, 85.80000000000001)
    ,@(28.82106135986731, 16.95999999999999, 42.62)
    ,@(9.369386401326686, 4.526666666666666, 15.94666666666667)
    ,@(89.39525704809283, 80.20000000000002, 94.55333333333333)
    ,@(2.454494195688217, 0.7466666666666669, 4.866666666666666)
    ,@(25.88673300165833, 13.86666666666666, 38.65999999999999)
    ,@(21.62202321724696, 11.62666666666667, 34.26666666666668)
    ,@(9.360232172470973, 4.133333333333335, 16.28666666666667)
    ,@(87.20245439469308, 77.89999999999999, 92.99333333333331)
    ,@(64.88325041459365, 50.39999999999998, 75.97333333333337)
    ,@(2.656981757877268, 0.9133333333333336, 4.993333333333333)
    ,@(49.74411276948576, 36.12666666666667, 64.74666666666664)
    ,@(83.27635157545595, 74.65333333333338, 90.7)
    ,@(86.53741293532333, 77.39333333333336, 92.67333333333333)
    ,@(83.82484245439458, 75.10666666666667, 90.40666666666671)
    ,@(68.52769485903805, 56.39333333333332, 79.15333333333332)
    ,@(26.83074626865665, 15.63333333333333, 39.40666666666667)
    ,@(63.17840796019896, 48.07333333333334, 76.26000000000001)
    ,@(10.51174129353234, 4.853333333333333, 17.46)
    ,@(83.69744610281919, 73.75999999999999, 90.28666666666665)
    ,@(52.35466003316746, 37.44666666666669, 65.32666666666665)
    ,@(7.852404643449415, 3.380000000000002, 14.4)
    ,@(61.19840796019894, 47.38, 74.38000000000001)
    ,@(77.93956882255384, 67.16000000000003, 86.42000000000002)
    ,@(49.74411276948576, 36.12666666666667, 64.74666666666664)
    ,@(74.11674958540617, 62.93333333333333, 83.92000000000002)
    ,@(78.37611940298497, 68.22666666666666, 86.57333333333332)
    ,@(82.73326699834145, 73.51333333333334, 90.56000000000003)
    ,@(89.01227197346593, 79.34, 94.72666666666665)
    ,@(51.51800995024865, 36.82, 64.99333333333333)
    ,@(7.153067993366496, 3.026666666666667, 12.37333333333333)
    ,@(10.48941956882255, 4.926666666666666, 17.48)
    ,@(1.354958540630182, 0.2400000000000001, 2.873333333333332)
    ,@(81.7224212271973, 72.56, 88.69333333333336)
    ,@(52.24165837479263, 36.86, 66.12666666666668)
    ,@(2.17641791044775, 0.7333333333333337, 4.206666666666669)
    ,@(9.363349917081253, 4.406666666666667, 15.82)
    ,@(83.56461028192361, 74.61333333333333, 90.34)
    ,@(86.48019900497505, 77.28, 92.66666666666664)
    ,@(10.59160862354891, 4.693333333333334, 17.93333333333334)
    ,@(76.50892205638459, 64.30666666666667, 85.74666666666664)
    ,@(2.494759535655047, 0.786666666666667, 4.973333333333334)
    ,@(83.84776119402974, 74.87333333333332, 90.4666666666667)
    ,@(77.29389718076278, 66.14000000000001, 86.23999999999999)
    ,@(84.02560530679926, 74.12000000000003, 90.86666666666665)
    ,@(2.501791044776108, 0.8066666666666669, 4.913333333333335)
    ,@(86.03618573797667, 76.84666666666666, 92.28666666666668)
    ,@(68.27260364842444, 56.42, 79.06666666666666)
    ,@(67.515655058043, 55.21333333333336, 78.63333333333335)
    ,@(8.995489220563837, 3.966666666666666, 15.31333333333333)
    ,@(89.33330016583741, 80.06666666666666, 94.61333333333336)
    ,@(7.641326699834153, 3.5, 13.69333333333333)
    ,@(40.56932006633486, 26.66666666666667, 55.50666666666668)
    ,@(46.42683250414577, 32.74666666666666, 61.66666666666665)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
    $ws.Cells.Item($row, 3).Value = $values[$i][2]
}
